$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 37 and row 38 data (Aptos <-> WEMIXTOKEN)
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.428"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.25%  "

# Update price/volume values for all other rows
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.943.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4953"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.86%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07228"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.054"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.924"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.736.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.831"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001032"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06372"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.006.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.049"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.930.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09450"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.575"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.98%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.359"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05917"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02186"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1990"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5983"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.106"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.444"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.58%  "

# Row 45
$ws.Range("E45").Value = "  -1.57%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.578"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5619"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.846"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06652"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "

# Row 51
$ws.Range("E51").Value = "  -4.90%  "
